$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content: swap in the new recipient e-mail addresses ---
$ws.Range("B2").Value = "NumilenMercado@yandex.com"
$ws.Range("B3").Value = "lupascugabrielcristian@gmail.com"

# --- Re-stamp the used range with the "Normal" cell style. Visually this
#     is a no-op (same font/fill/border/number format), but it is what
#     re-saving the workbook in Excel does: every cell picks up a freshly
#     written style record instead of the old implicit default. ---
$ws.Range("A1:C8").Style = "Normal"

# --- Widen column B so the longer e-mail addresses are fully visible;
#     column C is left essentially as-is. ---
$ws.Columns.Item(2).ColumnWidth = 42.96
$ws.Columns.Item(3).ColumnWidth = 21.38

# --- Selection / zoom left over from the editing session ---
[void]$ws.Range("B4").Select()
$excel.ActiveWindow.Zoom = 280
